# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K") values for rows 2-39 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 5
    3  = 10
    4  = 4
    5  = 8
    6  = 5
    7  = 6
    8  = 8
    9  = 6
    10 = 3
    11 = 2
    12 = 12
    13 = 3
    14 = 13
    15 = 2
    16 = 9
    17 = 8
    18 = 9
    19 = 9
    20 = 5
    21 = 3
    22 = 8
    23 = 8
    24 = 7
    25 = 8
    26 = 7
    27 = 5
    28 = 5
    29 = 6
    30 = 6
    31 = 6
    32 = 5
    33 = 10
    34 = 5
    35 = 4
    36 = 5
    37 = 5
    38 = 3
    39 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
